$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$ref, [string]$val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextCell "D2" "30.401.44"
Set-TextCell "E2" "  -0.76%  "
Set-TextCell "D3" "1.892.76"
Set-TextCell "E3" "  +0.24%  "
Set-TextCell "E4" "  -0.06%  "
Set-TextCell "D5" "237.85"
Set-TextCell "E5" "  +0.51%  "
Set-TextCell "E6" "  +0.01%  "
Set-TextCell "D7" "0.4894"
Set-TextCell "D8" "0.2931"
Set-TextCell "E8" "  +0.35%  "
Set-TextCell "D9" "0.06696"
Set-TextCell "E9" "  +0.33%  "
Set-TextCell "D10" "1.852.31"
Set-TextCell "E10" "  -2.01%  "
Set-TextCell "E11" "  +1.58%  "
Set-TextCell "D12" "0.07349"
Set-TextCell "E12" "  +1.50%  "
Set-TextCell "D13" "5.121"
Set-TextCell "E13" "  +2.35%  "
Set-TextCell "D14" "87.53"
Set-TextCell "E14" "  -1.91%  "
Set-TextCell "D15" "0.6627"
Set-TextCell "E15" "  -0.04%  "
Set-TextCell "D16" "30.373.30"
Set-TextCell "E16" "  -0.71%  "
Set-TextCell "D17" "13.41"
Set-TextCell "E17" "  +3.23%  "
Set-TextCell "D18" "0.000007815"
Set-TextCell "E18" "  -0.97%  "
Set-TextCell "E19" "  +0.08%  "
Set-TextCell "D20" "2.145.50"
Set-TextCell "E20" "  +0.47%  "
Set-TextCell "D21" "5.314"
Set-TextCell "E21" "  +11.97%  "
Set-TextCell "D22" "1.001"
Set-TextCell "E22" "  -0.08%  "
Set-TextCell "D23" "189.70"
Set-TextCell "E23" "  -1.23%  "
Set-TextCell "D24" "6.101"
Set-TextCell "E24" "  +0.19%  "
Set-TextCell "D25" "9.459"
Set-TextCell "E25" "  +1.69%  "
Set-TextCell "D26" "163.34"
Set-TextCell "E26" "  +2.20%  "
Set-TextCell "E27" "  -0.58%  "
Set-TextCell "D28" "1.926"
Set-TextCell "D29" "1.467"
Set-TextCell "E29" "  +4.53%  "
Set-TextCell "D30" "4.348"
Set-TextCell "E30" "  +2.15%  "
Set-TextCell "E31" "  +1.46%  "
Set-TextCell "D32" "4.031"
Set-TextCell "E32" "  +2.39%  "
Set-TextCell "D33" "0.05181"
Set-TextCell "E33" "  -0.18%  "
Set-TextCell "D34" "0.7387"
Set-TextCell "E34" "  +0.88%  "
Set-TextCell "E35" "  +1.13%  "
Set-TextCell "D36" "2.717"
Set-TextCell "E36" "  +1.19%  "
Set-TextCell "D37" "0.01809"
Set-TextCell "E37" "  -0.65%  "
Set-TextCell "D38" "2.666"
Set-TextCell "E38" "  -0.14%  "
Set-TextCell "D39" "0.9207"
Set-TextCell "E39" "  -0.35%  "
Set-TextCell "E40" "  -0.36%  "
Set-TextCell "D41" "0.4378"
Set-TextCell "E41" "  -0.62%  "
Set-TextCell "D42" "5.918"
Set-TextCell "E42" "  +3.19%  "
Set-TextCell "D43" "105.97"
Set-TextCell "E43" "  +1.60%  "
Set-TextCell "D44" "0.9920"
Set-TextCell "E44" "  -0.74%  "
Set-TextCell "D45" "68.23"
Set-TextCell "E45" "  +19.16%  "
Set-TextCell "D46" "0.1369"
Set-TextCell "E46" "  +2.52%  "
Set-TextCell "D47" "7.577"
Set-TextCell "E47" "  +3.22%  "
Set-TextCell "B48" "EnergySwap"
Set-TextCell "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D48" "8.986"
Set-TextCell "E48" "  +3.44%  "
Set-TextCell "B49" "Elrond"
Set-TextCell "C49" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextCell "D49" "34.83"
Set-TextCell "E49" "  +4.77%  "
Set-TextCell "D51" "0.3935"
Set-TextCell "E51" "  -4.21%  "
